$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D; existing D..K shift right to E..L
$ws.Columns.Item(4).Insert()

# Copy the number-format/style of the (now-shifted) neighboring column E into
# the new column D, but only for the rows that actually contain data cells in
# that block (skip the section-header rows 5,6,37,79 which have no D..K cells).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

Write-Host "done"
